# Auto-generated Excel COM-interop script to apply xlsx diff
# Fixes p-value columns that were corrupted (wrong magnitude) back to the
# correct small scientific-notation values, and restores cursor selections.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Table S3A gene untreated")
$ws2 = $wb.Worksheets.Item("Table S3B exon-level untreated")
$ws3 = $wb.Worksheets.Item("Table S3C exon_ratio untreated")

# --- Table S3A gene untreated: correct column D p-values ---
$ws1.Range("D2").Value = [double]"3.2531002975804001E-6"
$ws1.Range("D3").Value = [double]"4.00753942194377E-6"
$ws1.Range("D4").Value = [double]"7.51005497478065E-6"
$ws1.Range("D5").Value = [double]"8.5165563582613692E-6"
$ws1.Range("D6").Value = [double]"1.0994203781808001E-5"
$ws1.Range("D7").Value = [double]"1.1097234614159599E-5"
$ws1.Range("D8").Value = [double]"1.2623185255327599E-5"
$ws1.Range("D9").Value = [double]"1.9125053155182802E-6"
$ws1.Range("D10").Value = [double]"2.7164215945597301E-6"
$ws1.Range("D11").Value = [double]"2.9210806871861297E-5"
$ws1.Range("D12").Value = [double]"2.9626127132929496E-5"
$ws1.Range("D13").Value = [double]"2.9885943517500501E-5"
$ws1.Range("D14").Value = [double]"3.5318813850789301E-5"
$ws1.Range("D15").Value = [double]"4.4243398800200406E-5"
$ws1.Range("D16").Value = [double]"4.5671283783782704E-5"
$ws1.Range("D17").Value = [double]"4.6277852612289301E-5"
$ws1.Range("D18").Value = [double]"4.88573534182493E-5"
$ws1.Range("D19").Value = [double]"5.6824925431303709E-5"
$ws1.Range("D20").Value = [double]"5.9269895903545E-5"
$ws1.Range("D21").Value = [double]"6.4000000000000011E-5"

# --- Table S3B exon-level untreated: correct column E p-values ---
$ws2.Range("E2").Value = [double]"8.923685658329411E-7"
$ws2.Range("E3").Value = [double]"9.5551333548059716E-7"
$ws2.Range("E4").Value = [double]"1.5543830513614501E-6"
$ws2.Range("E5").Value = [double]"1.6741966518010702E-6"
$ws2.Range("E6").Value = [double]"2.6912121549030498E-6"
$ws2.Range("E7").Value = [double]"3.3289576405802602E-6"
$ws2.Range("E8").Value = [double]"4.7104855372021508E-6"
$ws2.Range("E9").Value = [double]"5.1930297774477501E-6"
$ws2.Range("E10").Value = [double]"5.8575664521445202E-6"
$ws2.Range("E11").Value = [double]"6.8768973331679104E-6"
$ws2.Range("E12").Value = [double]"7.7502586659097706E-6"
$ws2.Range("E13").Value = [double]"1.0048940857322599E-5"
$ws2.Range("E14").Value = [double]"1.0445198837066299E-5"
$ws2.Range("E15").Value = [double]"1.0499178217759598E-5"
$ws2.Range("E16").Value = [double]"1.1745736990565699E-5"
$ws2.Range("E17").Value = [double]"1.2077661177523499E-5"
$ws2.Range("E18").Value = [double]"1.3354872956483099E-5"
$ws2.Range("E19").Value = [double]"1.3770009687068302E-5"
$ws2.Range("E20").Value = [double]"1.41717704764734E-5"
$ws2.Range("E21").Value = [double]"1.4982134705149301E-5"

# --- Table S3C exon_ratio untreated: correct column H p-values ---
$ws3.Range("H2").Value = [double]"1.05E-7"
$ws3.Range("H3").Value = [double]"1.18E-7"
$ws3.Range("H4").Value = [double]"1.3300000000000001E-7"
$ws3.Range("H5").Value = [double]"1.36E-7"
$ws3.Range("H6").Value = [double]"1.36E-7"
$ws3.Range("H7").Value = [double]"1.4000000000000001E-7"
$ws3.Range("H8").Value = [double]"1.4399999999999999E-7"
$ws3.Range("H9").Value = [double]"1.4399999999999999E-7"
$ws3.Range("H10").Value = [double]"1.4399999999999999E-7"
$ws3.Range("H11").Value = [double]"1.67E-7"
$ws3.Range("H12").Value = [double]"1.8800000000000002E-7"
$ws3.Range("H13").Value = [double]"2.8900000000000001E-7"
$ws3.Range("H14").Value = [double]"2.9500000000000003E-7"
$ws3.Range("H15").Value = [double]"3.3800000000000004E-7"
$ws3.Range("H16").Value = [double]"3.4700000000000002E-7"
$ws3.Range("H17").Value = [double]"3.7600000000000003E-7"
$ws3.Range("H18").Value = [double]"3.9200000000000002E-7"
$ws3.Range("H19").Value = [double]"4.0400000000000002E-7"
$ws3.Range("H20").Value = [double]"4.1162440899999998E-7"
$ws3.Range("H21").Value = [double]"4.1167600800000002E-7"

# --- Restore per-sheet cell selections (order matters: last Activate
#     call wins for which sheet stays tabSelected) ---
$null = $ws2.Range("F1:F1048576").Select()
$null = $ws3.Range("A2").Select()
$null = $ws1.Activate()
$null = $ws1.Range("D3").Select()
